$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'36.426.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.43%  "

$ws.Range("D3").Value = "'2.010.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.33%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "'251.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.10%  "

$ws.Range("D6").Value = "'0.638"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.28%  "

$ws.Range("D7").Value = "'61.44"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +12.69%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").Value = "'58.70"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.85%  "

$ws.Range("B10").Value = "Cardano"
$ws.Range("C10").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D10").Value = "'0.370"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.96%  "

$ws.Range("E11").Value = "  +0.63%  "

$ws.Range("D12").Value = "'0.103"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.22%  "

$ws.Range("D13").Value = "'0.895"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.18%  "

$ws.Range("E14").Value = "  +4.29%  "

$ws.Range("D15").Value = "'2.303.92"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.50%  "

$ws.Range("D16").Value = "'20.31"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +16.36%  "

$ws.Range("D17").Value = "'5.43"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.97%  "

$ws.Range("D18").Value = "'2.007.85"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.84%  "

$ws.Range("D19").Value = "'36.389.48"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.59%  "

$ws.Range("D20").Value = "'71.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.87%  "

$ws.Range("D21").Value = "'0.0₃0862"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.07%  "

$ws.Range("D22").Value = "'5.25"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.53%  "

$ws.Range("D23").Value = "'234.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.71%  "

$ws.Range("D24").Value = "'2.74"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +20.49%  "

$ws.Range("E25").Value = "  -0.02%  "

$ws.Range("E26").Value = "  -0.91%  "

$ws.Range("D27").Value = "'9.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.59%  "

$ws.Range("D28").Value = "'163.73"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.26%  "

$ws.Range("D29").Value = "'19.60"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.31%  "

$ws.Range("E30").Value = "  -1.06%  "

$ws.Range("D31").Value = "'5.10"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.23%  "

$ws.Range("D32").Value = "'0.111"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +23.28%  "

$ws.Range("E33").Value = "  +0.82%  "

$ws.Range("E34").Value = "  +5.87%  "

$ws.Range("D35").Value = "'0.0607"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.80%  "

$ws.Range("E36").Value = "  +10.96%  "

$ws.Range("E37").Value = "  -0.04%  "

$ws.Range("D38").Value = "'1.80"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.20%  "

$ws.Range("E39").Value = "  +17.66%  "

$ws.Range("E40").Value = "  +13.89%  "

$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "'2.78"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +23.47%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'1.23"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.69%  "

$ws.Range("E44").Value = "  +2.95%  "

$ws.Range("D45").Value = "'0.0215"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.58%  "

$ws.Range("D46").Value = "'7.99"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.57%  "

$ws.Range("D47").Value = "'16.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.40%  "

$ws.Range("D48").Value = "'94.23"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.46%  "

$ws.Range("D49").Value = "'1.423.15"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.78%  "

$ws.Range("D50").Value = "'2.91"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.99%  "

$ws.Range("D51").Value = "'46.90"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.58%  "
